$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-6 from 45224 to 45233
$ws.Range("C2:C6").Value = 45233
